$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1101.71969107791
$ws.Range("C3").Value = 16.50033441817651
$ws.Range("C4").Value = 9981640.870446853
$ws.Range("C5").Value = 5098.905377874559
$ws.Range("C6").Value = 167.7131847245725
$ws.Range("C7").Value = 26911137.68433562
$ws.Range("C8").Value = 1.988366255557708
$ws.Range("C10").Value = 9570.795881577404
$ws.Range("C13").Value = 9476.819960762952
$ws.Range("C14").Value = 0.7741940793718595
$ws.Range("C15").Value = 0.01942137731683652
$ws.Range("C16").Value = 1.030935261260624
$ws.Range("C17").Value = 0.01739538833450718
$ws.Range("C18").Value = 35990.15460923022
$ws.Range("C19").Value = 650.2645868542827
$ws.Range("C20").Value = 350565217.0365177
$ws.Range("C21").Value = 175214.8226370455
$ws.Range("C22").Value = 80808.30563198758
$ws.Range("C23").Value = 30781676152.91018
$ws.Range("C24").Value = 0.7275316865904939
$ws.Range("C26").Value = 3785.123608469989
$ws.Range("C27").Value = 0.7100226467075748
$ws.Range("C28").Value = 0.01304453275894353
$ws.Range("C29").Value = 3748.585088812148
$ws.Range("C30").Value = 0.8186384684541966
$ws.Range("C31").Value = 0.01165619169488945
$ws.Range("C32").Value = 0.6734476128704446
$ws.Range("C33").Value = 0.01043847243578569
$ws.Range("C35").Value = 0.06652638069348274
$ws.Range("C36").Value = 6675.764077720934
$ws.Range("C37").Value = 2.960025942834124
$ws.Range("C38").Value = 0.06002665629918163
$ws.Range("C40").Value = 0.5967648023438751
$ws.Range("C41").Value = 0.02946571340099468
$ws.Range("C42").Value = 4863.66971617547
$ws.Range("C46").Value = 0.01368964270691804
$ws.Range("C47").Value = 14.36727565378699
$ws.Range("C48").Value = 9.200681933715055
$ws.Range("C49").Value = 0.01047853129480764
$ws.Range("C50").Value = 0.03537336823135129
$ws.Range("C51").Value = 26.20830182689003
$ws.Range("C52").Value = 8.670413553089622
$ws.Range("C53").Value = 0.01792059537637969
$ws.Range("C54").Value = 48227038.72102559
$ws.Range("C55").Value = 341883903.0945713
$ws.Range("C56").Value = 0.01982061956315294
$ws.Range("C57").Value = 87.45989435086626
$ws.Range("C58").Value = 213.6127448828788
$ws.Range("C59").Value = 0.01699975027811759
$ws.Range("C60").Value = 0.02911208515149327
$ws.Range("C61").Value = 26.21342208640771
$ws.Range("C62").Value = 10.38498293237319
$ws.Range("C63").Value = 0.04176107288800473
$ws.Range("C64").Value = 46562633.90746562
$ws.Range("C65").Value = 212270025.108244
$ws.Range("C66").Value = 0.5148039188497781
$ws.Range("C67").Value = 0.5149910913365231
$ws.Range("C68").Value = 0.01077003109689171
$ws.Range("C69").Value = 0.01061455877641167
$ws.Range("C74").Value = 4276.275497911609
$ws.Range("C75").Value = 4315.300441781877
$ws.Range("C76").Value = 0.8730603803542444
$ws.Range("C77").Value = 0.8014520792494723
$ws.Range("C78").Value = 0.01650507044530855
$ws.Range("C79").Value = 0.01735239056592375
$ws.Range("C81").Value = 0.7974725467761707
$ws.Range("C84").Value = 7890.722950254709
$ws.Range("C85").Value = 8185.427444549284
$ws.Range("C86").Value = 0.5243123006247158
$ws.Range("C87").Value = 0.532365290780599
$ws.Range("C90").Value = 0.6564044780027963
$ws.Range("C91").Value = 0.5229491375241281
$ws.Range("C94").Value = 1394148.530056586
$ws.Range("C95").Value = 1242207.864805819
$ws.Range("C96").Value = 0.8639314804057794
$ws.Range("C97").Value = 0.6938775654771154
$ws.Range("C100").Value = 0.8242527216766448
$ws.Range("C104").Value = 809884.2898824479
$ws.Range("C105").Value = 1432573.3096264
$ws.Range("C106").Value = 0.8639290053995139
$ws.Range("C107").Value = 0.6696386831182918
$ws.Range("C110").Value = 0.8242526945273037
$ws.Range("C114").Value = 7620.972452829933
$ws.Range("C115").Value = 7028.51618514347
$ws.Range("C116").Value = 0.9415811887921495
$ws.Range("C117").Value = 0.01535077065020322
$ws.Range("C118").Value = 0.9850998564439555
$ws.Range("C119").Value = 0.01911127630413923
$ws.Range("C120").Value = 1.077352350037336
$ws.Range("C121").Value = 0.02213519130350525
$ws.Range("C122").Value = 5079.146816273959
$ws.Range("C123").Value = 0.6329081018208389
$ws.Range("C124").Value = 0.0144537314725022
$ws.Range("C125").Value = 6676.846749475254
$ws.Range("C126").Value = 0.8625727873660767
$ws.Range("C127").Value = 0.01718549974873273
$ws.Range("C128").Value = 6344.231895294703
$ws.Range("C129").Value = 0.8897060289691914
$ws.Range("C130").Value = 0.01617188602269905
$ws.Range("C131").Value = 6501.158311649955
$ws.Range("C132").Value = 0.9162146816014356
$ws.Range("C133").Value = 0.02234884021828088
$ws.Range("C134").Value = 1.220338303785464
$ws.Range("C135").Value = 0.02447613228335998
$ws.Range("C136").Value = 1.33079603528621
$ws.Range("C137").Value = 0.01800666967449917
$ws.Range("C138").Value = 7348.261531580207
$ws.Range("C139").Value = 1.093564845650243
$ws.Range("C140").Value = 0.02534778834317109
$ws.Range("C141").Value = 7479.059057970027
$ws.Range("C142").Value = 1.077161194326148
$ws.Range("C143").Value = 0.01947425138547508
$ws.Range("C144").Value = 8318.398976607101
$ws.Range("C145").Value = 1.10995726851603
$ws.Range("C146").Value = 0.02062317677837518
$ws.Range("C147").Value = 8132.376455300657
